# Powerpoint writer: consolidate text run nodes.
# Merge the "First"/" " runs (and "Third"/" " runs) in the title
# placeholders of slide 1 and slide 3 into a single run each, e.g.
# "First" + " " + "slide"  ->  "First " + "slide"
# "Third" + " " + "slide"  ->  "Third " + "slide"
# This reduces the run count without changing the visible text.

$p = $ppt.ActivePresentation

$targets = @(
    @{ Slide = 1; Word = "First" },
    @{ Slide = 3; Word = "Third" }
)

foreach ($target in $targets) {
    $s = $p.Slides.Item($target.Slide)
    $shape = $s.Shapes.Item(1)
    $tr = $shape.TextFrame.TextRange

    # Re-write the "<Word> " prefix (the first two runs: "<Word>" and " ")
    # as a single run of text, leaving the trailing "slide" run untouched.
    $prefixLen = $target.Word.Length + 1
    $prefix = $tr.Characters(1, $prefixLen)
    $prefix.Text = $target.Word + " "
}
